$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.325.17'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '1.870.85'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '236.24'
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '0.4708'
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").Value = '0.2900'
$ws.Range("E8").Value = '  +2.06%  '

$ws.Range("D9").Value = '0.06628'
$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").Value = '21.75'
$ws.Range("E10").Value = '  -0.49%  '

$ws.Range("D11").Value = '0.08015'
$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("D12").Value = '97.43'
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").Value = '1.878.65'
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("D14").Value = '5.160'
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = '0.6883'
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").Value = '274.26'
$ws.Range("E16").Value = '  -1.99%  '

$ws.Range("D17").Value = '30.311.86'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").Value = '14.07'
$ws.Range("E18").Value = '  +6.52%  '

$ws.Range("D19").Value = '0.000007724'
$ws.Range("E19").Value = '  +5.64%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").Value = '2.122.82'
$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("D22").Value = '5.321'
$ws.Range("E22").Value = '  -1.77%  '

$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("D24").Value = '6.227'
$ws.Range("E24").Value = '  +1.18%  '

$ws.Range("D25").Value = '167.69'
$ws.Range("E25").Value = '  +0.91%  '

$ws.Range("D26").Value = '9.285'
$ws.Range("E26").Value = '  +1.30%  '

$ws.Range("D27").Value = '19.01'
$ws.Range("E27").Value = '  -0.48%  '

$ws.Range("D28").Value = '1.964'
$ws.Range("E28").Value = '  +1.39%  '

$ws.Range("D29").Value = '1.375'
$ws.Range("E29").Value = '  -1.03%  '

$ws.Range("D30").Value = '0.09950'
$ws.Range("E30").Value = '  +2.18%  '

$ws.Range("D31").Value = '4.373'
$ws.Range("E31").Value = '  -0.73%  '

$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").Value = '4.094'
$ws.Range("E33").Value = '  +0.10%  '

$ws.Range("D34").Value = '0.04712'
$ws.Range("E34").Value = '  -0.42%  '

$ws.Range("D35").Value = '1.135'
$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").Value = '0.7038'
$ws.Range("E36").Value = '  -0.73%  '

$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  -0.26%  '

$ws.Range("D38").Value = '0.01883'
$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("D39").Value = '2.633'
$ws.Range("E39").Value = '  +2.19%  '

$ws.Range("D40").Value = '6.326'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").Value = '73.49'
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("D42").Value = '1.964'
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4172'
$ws.Range("E43").Value = '  -0.51%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.8419'
$ws.Range("E44").Value = '  -1.06%  '

$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").Value = '  -0.13%  '

$ws.Range("D46").Value = '103.75'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.114'
$ws.Range("E47").Value = '  -1.27%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.216'
$ws.Range("E48").Value = '  -2.17%  '

$ws.Range("D49").Value = '935.15'
$ws.Range("E49").Value = '  -3.80%  '

$ws.Range("D50").Value = '34.51'
$ws.Range("E50").Value = '  +1.06%  '

$ws.Range("E51").Value = '  +0.44%  '
